$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for rows 3, 4, 5 in columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")
$rows = @(3, 4, 5)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Cyclic shift: new row3 = old row4, new row4 = old row5, new row5 = old row3
$mapping = @{ 3 = 4; 4 = 5; 5 = 3 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $orig[$src][$c]
    }
}
